# "Added Boss Health Bar"
#
# The document contains two "Add AI" bullet items (one highlighted red,
# under a different list, and one highlighted yellow, under the boss's
# AI-behavior list). Only the yellow-highlighted "Add AI" bullet (the one
# next to "Patrol Air" / "Shoot Bullets" / boss health-bar notes) needs to
# change its highlight color from yellow to (bright) green - both on the
# run of text itself and on the paragraph mark, matching how Word applies
# highlighting when the whole line (incl. the pilcrow) is selected.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Add AI"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0  # wdFindStop - don't wrap around, we walk all matches once

$found = $rng.Find.Execute()
while ($found) {
    # wdYellow highlight reads back as the hex string "#FFFF00"
    if ($rng.Font.HighlightColorIndex -eq "#FFFF00") {
        # Expand to the full paragraph (including the paragraph mark) so the
        # highlight change applies to both the run and the paragraph mark's
        # run properties, exactly like selecting the whole line in Word.
        $paraRng = $rng.Paragraphs(1).Range
        $paraRng.Font.HighlightColorIndex = 4  # wdBrightGreen -> w:val="green"
    }
    $found = $rng.Find.Execute()
}
